# Apply the commit's edit: reset the "Created in sprint" (column B) values
# for a block of story rows back to 0, and move the sheet selection to B91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(32,33,34,36,38,40,41,42,43,44,46,48,50,51,52,54,56,57,58,60,61,62,63,64,65,67,68,70,71,72,74,76,77,78,79,81,82,83,85,86,87,88,89,91)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = 0
}

# Move the active selection to B91 (matches the saved view state in the diff).
$ws.Activate() | Out-Null
$ws.Range("B91").Select() | Out-Null
